# Apply the edits described in the commit:
#  - rename the 4th sheet ("工作表4") to "day4"
#  - add a new "当天总分" (daily total score) column header to the end of
#    row 1 on every sheet

$wb = $excel.ActiveWorkbook

# 1) Rename the 4th worksheet tab.
$daySheet = $wb.Worksheets.Item(4)
$daySheet.Name = "day4"

# 2) Add the new "当天总分" header cell after the last used column of row 1
#    on every worksheet.
foreach ($ws in $wb.Worksheets) {
    $lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End(-4159).Column
    $ws.Cells.Item(1, $lastCol + 1).Value = "当天总分"
}
